# Update countries & provincias Spain
# Updates case numbers for several countries and swaps the row order of
# "Mali" / "El Salvador" (El Salvador now has more cases, so it moves
# above Mali) while refreshing El Salvador's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Kazajistan (row 67): Casos activos, Recuperados
$ws.Range("D67").Value = 285
$ws.Range("E67").Value = 1168

# Uzbekistan (row 71): Casos totales, Nuevos casos, Recuperados
$ws.Range("B71").Value = 1380
$ws.Range("C71").Value = 31
$ws.Range("E71").Value = 1247

# Lituania (row 76): Casos totales, Nuevos casos, Casos activos, Recuperados
$ws.Range("B76").Value = 1149
$ws.Range("C76").Value = 21
$ws.Range("D76").Value = 210
$ws.Range("E76").Value = 907

# Taiwan (row 107): Casos activos, Recuperados
$ws.Range("D107").Value = 166
$ws.Range("E107").Value = 223

# Georgia (row 110): Casos totales, Nuevos casos, Casos activos, Recuperados
$ws.Range("B110").Value = 370
$ws.Range("C110").Value = 22
$ws.Range("D110").Value = 77
$ws.Range("E110").Value = 290

# Row 124 now becomes "El Salvador" (was Mali) with updated figures
$ws.Range("A124").Value = "El Salvador"
$ws.Range("B124").Value = 177
$ws.Range("C124").Value = 13
$ws.Range("D124").Value = 33
$ws.Range("E124").Value = 137
$ws.Range("F124").Value = 2
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 7

# Row 125 now becomes "Mali" (was El Salvador) with its original figures
$ws.Range("A125").Value = "Mali"
$ws.Range("B125").Value = 171
$ws.Range("C125").Value = 0
$ws.Range("D125").Value = 34
$ws.Range("E125").Value = 124
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 13
